# Weekly fruit/vegetable data update.
# A new daily observation is inserted as row 73 (pushing the existing rows
# 73-139 down to 74-140, growing the used range from A1:R139 to A1:R140).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 73; everything below shifts down by one.
$ws.Rows("73:73").Insert()

# Populate the newly inserted row 73 with the new observation.
$ws.Range("A73").Value = 6
$ws.Range("B73").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C73").Value = "Metropolitana"
$ws.Range("D73").Value = 44539
$ws.Range("E73").Value = 13
$ws.Range("F73").Value = 100112001
$ws.Range("G73").Value = "Berenjena"
$ws.Range("H73").Value = "Sin especificar"
$ws.Range("I73").Value = "Primera"
$ws.Range("J73").Value = 180
$ws.Range("K73").Value = 9000
$ws.Range("L73").Value = 10000
$ws.Range("M73").Value = 9556
$ws.Range("N73").Value = "`$/caja 50 unidades"
$ws.Range("O73").Value = "Región de Arica y Parinacota"
$ws.Range("P73").Value = 191
$ws.Range("Q73").Value = 50
$ws.Range("R73").Value = "Hortaliza"
